# Sprint 5 closeout: trim the "Datos" sample data down to two scenarios
# (Acierto / Alterno) and drop the now-unused transactional detail columns
# for the remaining rows, matching the finalized data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Remove the third data row (former row 4: "Alterno" / CLAVE BLOQUEADA case).
# This shifts nothing below it, shrinks the used range to A1:J3, and is
# reflected automatically in the sheet dimension.
$ws.Rows.Item(4).Delete()

# The remaining two data rows only keep their idCaso/orientacion values;
# the rest of the transactional columns (codigoTransaccion ... tipoDocumento)
# are cleared out for both rows.
$ws.Range("C2:J3").ClearContents()
